# Adds an "Estadísticas" (Statistics) section right after the
# "Dispositivo" / "Sensores" summary block that ends with "Humedad (%)".
#
# The new paragraph mirrors the look of the existing "Dispositivo" and
# "Sensores" paragraphs already in the document: a bold 13pt title run
# followed by five 12pt explanatory lines, each pair separated by a
# manual line break (<w:br/>) inside a single paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that ends with "Humedad (%)" (the last line of
# the "Sensores" block) without hard-coding a paragraph index.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Humedad (%)", $true, $false, $false,
                                    $false, $false, $true, 1, $false,
                                    "", 0)
if (-not $found) {
    throw "Could not find anchor text 'Humedad (%)'"
}

$targetParaIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -le $searchRange.Start -and `
        $candidate.Range.End -ge $searchRange.End) {
        $targetParaIndex = $i
        break
    }
}
if ($targetParaIndex -eq 0) {
    throw "Could not resolve the paragraph containing 'Humedad (%)'"
}

$anchorPara = $d.Paragraphs.Item($targetParaIndex)

# Insert a brand-new paragraph right after it; Word seeds its pPr from
# the paragraph it was split from, which is exactly the spacing/
# indentation this new block should carry.
$anchorPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetParaIndex + 1)
$newRange = $newPara.Range

# Build the paragraph body (pPr + runs) as literal WordprocessingML and
# drop it in with InsertXML so every run gets exact run-formatting
# (including complex-script siblings bCs/szCs) that simple Font.*
# properties can't reach.
$xml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind/><w:spacing w:before="200" w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Estad&#237;sticas</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:color w:val="002060"/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Media: Media de los valores en el an&#225;lisis</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:color w:val="002060"/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Max: Valor m&#225;ximo horario</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:color w:val="002060"/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Min: Valor m&#237;nimo horario</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:color w:val="002060"/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Fecha Max: Fecha cuando se report&#243; el valor m&#225;ximo horario</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:color w:val="002060"/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Fecha Min: Fecha cuando se report&#243; el valor m&#237;nimo horario</w:t></w:r></w:p></w:body></w:document>
'@

$newRange.InsertXML($xml)

Write-Output "Inserted Estadisticas section after paragraph $targetParaIndex"
